# Retraining the models for GESS
# Update Consumption_Actual data: new values for rows 2-47, plus two new
# rows (48-49) appended, extending the dimension from A1:B47 to A1:B49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6338,46048.95833333334),
    @(6316,46048.96875),
    @(6254,46048.97916666666),
    @(6184,46048.98958333334),
    @(6203,46049),
    @(6097,46049.01041666666),
    @(6051,46049.02083333334),
    @(6023,46049.03125),
    @(5920,46049.04166666666),
    @(5891,46049.05208333334),
    @(5890,46049.0625),
    @(5895,46049.07291666666),
    @(5978,46049.08333333334),
    @(6016,46049.09375),
    @(5981,46049.10416666666),
    @(5963,46049.11458333334),
    @(5992,46049.125),
    @(5987,46049.13541666666),
    @(6041,46049.14583333334),
    @(6082,46049.15625),
    @(6183,46049.16666666666),
    @(6225,46049.17708333334),
    @(6297,46049.1875),
    @(6441,46049.19791666666),
    @(6580,46049.20833333334),
    @(6773,46049.21875),
    @(6875,46049.22916666666),
    @(7119,46049.23958333334),
    @(7472,46049.25),
    @(7666,46049.26041666666),
    @(7801,46049.27083333334),
    @(7872,46049.28125),
    @(8073,46049.29166666666),
    @(8266,46049.30208333334),
    @(8344,46049.3125),
    @(8343,46049.32291666666),
    @(8515,46049.33333333334),
    @(8550,46049.34375),
    @(8570,46049.35416666666),
    @(8577,46049.36458333334),
    @(8558,46049.375),
    @(8567,46049.38541666666),
    @(8547,46049.39583333334),
    @(8510,46049.40625),
    @(8447,46049.41666666666),
    @(8410,46049.42708333334),
    @(8464,46049.4375),
    @(8385,46049.44791666666)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $data[$i][1]
    # Column B cells already use numFmtId 164 ("YYYY-MM-DD HH:MM:SS"); make
    # sure the two newly appended rows (48-49) get the same number format.
    $cellB.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
